{"js": "// Remove the \"pa\u00eds\" (country) placeholder from each row of the\n// \"Formaci\u00f3n acad\u00e9mica\" (education) table. Each row's \"Instituci\u00f3n y pa\u00eds\"\n// cell originally reads \"${insX} ${paisX}\"; after the edit it should read\n// just \"${insX}\" \u2014 i.e. drop the trailing \" ${paisX}\" run sequence.\nconst body = context.document.body;\n\nconst suffixes = [\"L\", \"E\", \"M\", \"D\"];\n\nfor (const suffix of suffixes) {\n  const needle = \" ${pais\" + suffix + \"}\";\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"pa\u00eds\" (country) placeholder from each row of the\n# \"Formaci\u00f3n acad\u00e9mica\" (education) table. Each row's \"Instituci\u00f3n y pa\u00eds\"\n# cell originally reads \"${insX} ${paisX}\"; after the edit it should read\n# just \"${insX}\" \u2014 i.e. drop the trailing \" ${paisX}\" run sequence.\n$d = $word.ActiveDocument\n\n$suffixes = @('L', 'E', 'M', 'D')\n\nforeach ($suffix in $suffixes) {\n    $needle = ' ${pais' + $suffix + '}'\n    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, '', 2) | Out-Null\n}\n"}
